$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("sku", "name", "quantity", "cost_per", "total_cost")

for ($row = 2; $row -le 4; $row++) {
    for ($col = 1; $col -le 5; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
